# Weekly update for Fruta/Cereza sheet:
# Insert 7 new rows before row 433 (pushing the old 433-444 block down to
# 440-451) and populate the new rows 433-439 with the latest week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 7 blank rows above the current row 433 - shifts rows 433:444 -> 440:451
$ws.Range("A433:A439").EntireRow.Insert()

# Common values shared by every row in this block
$mercadoId = 8
$mercado   = "Terminal La Palmera de La Serena"
$region    = "Coquimbo"
$codreg    = 4
$tipo      = "Fruta"
$productoId = 100103
$producto  = "Frutos de hueso (carozo)"
$categoriaId = 100103001
$categoria = "Cereza"

# New rows' variable data: D(fecha), K(variedad), L(calidad), M(volumen),
# N(precio minimo), O(precio maximo), P(precio promedio ponderado),
# Q(unidad comercializacion), R(origen), S(precio $/Kg), T(kg/unidad)
$newRows = @(
    @(44931, "Brooks",      "Especial", 400, 12000, 13000, 12500, "`$/caja 15 kilos", "Provincia de Curicó", 833, 15),
    @(44931, "Brooks",      "Primera",  400, 10000, 11000, 10500, "`$/caja 15 kilos", "Provincia de Curicó", 700, 15),
    @(44931, "Lapins",      "Especial", 240, 13000, 14000, 13500, "`$/caja 15 kilos", "Provincia de Curicó", 900, 15),
    @(44931, "Lapins",      "Primera",  400, 11000, 12000, 11500, "`$/caja 15 kilos", "Provincia de Curicó", 767, 15),
    @(44931, "Lapins",      "Segunda",  360, 9000,  10000, 9500,  "`$/caja 15 kilos", "Provincia de Curicó", 633, 15),
    @(44931, "Sweet Heart", "Primera",  400, 11000, 12000, 11500, "`$/caja 15 kilos", "Provincia de Curicó", 767, 15),
    @(44931, "Sweet Heart", "Segunda",  360, 9000,  10000, 9500,  "`$/caja 15 kilos", "Provincia de Curicó", 633, 15)
)

$r = 433
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value  = $mercadoId
    $ws.Cells.Item($r, 2).Value  = $mercado
    $ws.Cells.Item($r, 3).Value  = $region
    $ws.Cells.Item($r, 4).Value  = $row[0]
    $ws.Cells.Item($r, 5).Value  = $codreg
    $ws.Cells.Item($r, 6).Value  = $tipo
    $ws.Cells.Item($r, 7).Value  = $productoId
    $ws.Cells.Item($r, 8).Value  = $producto
    $ws.Cells.Item($r, 9).Value  = $categoriaId
    $ws.Cells.Item($r, 10).Value = $categoria
    $ws.Cells.Item($r, 11).Value = $row[1]
    $ws.Cells.Item($r, 12).Value = $row[2]
    $ws.Cells.Item($r, 13).Value = $row[3]
    $ws.Cells.Item($r, 14).Value = $row[4]
    $ws.Cells.Item($r, 15).Value = $row[5]
    $ws.Cells.Item($r, 16).Value = $row[6]
    $ws.Cells.Item($r, 17).Value = $row[7]
    $ws.Cells.Item($r, 18).Value = $row[8]
    $ws.Cells.Item($r, 19).Value = $row[9]
    $ws.Cells.Item($r, 20).Value = $row[10]
    $r++
}
